# Fruta / hortaliza, semanal
# Re-applies the weekly refresh of the "Tuna" price table: rows 2-7 are
# updated in place so each row now carries the data that, in the previous
# snapshot, lived in a different row (the underlying source table was
# re-pulled and the local rows got overwritten with newer/rotated records).
# Only the columns that actually differ between snapshots are touched:
# D (Fecha), L (Calidad), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado), R (Origen), S (Precio $/Kg).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{ D = 44250; L = "Primera"; M = 200; N = 14000; O = 15000; P = 14500; R = "Región Metropolitana";    S = 806 }
    3 = @{ D = 45072; L = "Segunda"; M = 100; N = 16000; O = 16000; P = 16000; R = "Provincia de Chacabuco";  S = 889 }
    4 = @{ D = 45072; L = "Segunda"; M = 100; N = 17000; O = 17000; P = 17000; R = "Provincia de Limarí";     S = 944 }
    5 = @{ D = 44257; L = "Primera"; M = 100; N = 14000; O = 15000; P = 14500; R = "Región Metropolitana";    S = 806 }
    6 = @{ D = 44252; L = "Primera"; M = 120; N = 13000; O = 14000; P = 13500; R = "Región Metropolitana";    S = 750 }
    7 = @{ D = 44253; L = "Primera"; M = 160; N = 14000; O = 15000; P = 14500; R = "Región Metropolitana";    S = 806 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("R$row").Value = $vals.R
    $ws.Range("S$row").Value = $vals.S
}
